$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheet1): F2 1050 -> 1052, F3 21 -> 22
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1052
$ws1.Range("F3").Value = 22

# Update "全部类型" sheet (sheet4): F2 1050 -> 1052, F3 21 -> 22
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1052
$ws4.Range("F3").Value = 22
